# Rename severity-level headers (E1:L1 on every sheet) and drop the
# "Category" column (M) which only ever held the constant "primary".
$wb = $excel.ActiveWorkbook

$headerRenames = @{
    "E1" = "% severity levels 1-2"
    "F1" = "# severity levels 1-2"
    "G1" = "% severity level 3"
    "H1" = "# severity level 3"
    "I1" = "% severity level 4"
    "J1" = "# severity level 4"
    "K1" = "% severity level 5"
    "L1" = "# severity level 5"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in $headerRenames.Keys) {
        $ws.Range($addr).Value = $headerRenames[$addr]
    }
    $ws.Columns("M").Delete()
}
